$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation for Jengibre (Mercado Mayorista Lo Valledor de
# Santiago) needs to be inserted as row 125, pushing the existing rows
# 125-137 down to 126-138 (same as the data in the rest of the sheet just
# shifting to make room for the newest record at the top of this block).
$ws.Rows("125:125").Insert()

# Populate the newly-inserted row 125 with the new observation. The
# "dimension" columns (market id/name/region/category/origin/classification)
# match every other row for this product; only the date + price/volume
# figures differ.
$ws.Cells.Item(125, 1).Value = 6
$ws.Cells.Item(125, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(125, 3).Value = "Metropolitana"
$ws.Cells.Item(125, 4).Value = 45124
$ws.Cells.Item(125, 5).Value = 13
$ws.Cells.Item(125, 6).Value = 100114007
$ws.Cells.Item(125, 7).Value = "Jengibre"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 172
$ws.Cells.Item(125, 11).Value = 14000
$ws.Cells.Item(125, 12).Value = 15000
$ws.Cells.Item(125, 13).Value = 14552
$ws.Cells.Item(125, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(125, 15).Value = "Perú"
$ws.Cells.Item(125, 16).Value = 970
$ws.Cells.Item(125, 17).Value = 15
$ws.Cells.Item(125, 18).Value = "Hortaliza"
